# Applies the "Add files via upload" edit:
#   1. The subtitle "Projekt „2048“ von Georg Römmling und Florian Mansfeld"
#      had been split across multiple runs with spell-check proofErr markers
#      around "Römmling"; the new version is a single clean run (no proofErr).
#   2. The Score paragraph gains a clause explaining why the Highscore is
#      kept: ", um den Highscore mehrerer Spieler vergleichen zu können"
#      inserted right before the paragraph's closing period.

$d = $word.ActiveDocument

# --- 1. Subtitle: collapse the proofErr-split runs into one clean run ---
$quoteOpen  = [char]8222   # „
$quoteClose = [char]8220   # "
$subtitle   = "Projekt $quoteOpen" + "2048" + "$quoteClose von Georg Römmling und Florian Mansfeld"

$found1 = $d.Content.Find.Execute(
    $subtitle, $true, $false, $false, $false, $false, $true, 1, $false,
    $subtitle, 2
)
if (-not $found1) {
    throw "Could not find the subtitle text to normalize."
}

# --- 2. Score paragraph: insert the clause right before the final period ---
$scoreRange = $d.Content
$found2 = $scoreRange.Find.Execute(
    "Highscore-Ermittlung.", $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0
)
if (-not $found2) {
    throw "Could not find the Highscore-Ermittlung sentence."
}

# $scoreRange now spans the matched "Highscore-Ermittlung." text; collapse it
# to a point right before the trailing period and type the new clause there,
# the same way a cursor placed just before the "." would.
$scoreRange.Start = $scoreRange.End - 1
$scoreRange.Collapse(1)
$scoreRange.InsertBefore(", um den Highscore mehrerer Spieler vergleichen zu können")

# --- 3. Best-effort: register the built-in "comment" / "balloon text" style
#     family (Kommentarzeichen, Kommentartext(+Zchn), Kommentarthema(+Zchn),
#     Sprechblasentext(+Zchn)) that a resave of this document in a newer
#     Word pulls into styles.xml from the built-in style gallery. Wrapped in
#     try/catch so the (more important) content edits above always survive
#     even if a given style property isn't supported by this host.
try {
    $s1 = $d.Styles.Add("Kommentarzeichen", 2)
    $s1.NameLocal = "annotation reference"
    $s1.BaseStyle = "Absatz-Standardschriftart"
    $s1.Priority = 99
    $s1.UnhideWhenUsed = $true
    $s1.Font.Size = 8
    $s1.Font.SizeBi = 8

    $s2 = $d.Styles.Add("Kommentartext", 1)
    $s2.NameLocal = "annotation text"
    $s2.BaseStyle = "Standard"
    $s2.Priority = 99
    $s2.UnhideWhenUsed = $true
    $s2.LinkStyle = "KommentartextZchn"
    $s2.ParagraphFormat.LineSpacingRule = 0
    $s2.Font.Size = 10
    $s2.Font.SizeBi = 10

    $s3 = $d.Styles.Add("KommentartextZchn", 2)
    $s3.NameLocal = "Kommentartext Zchn"
    $s3.BaseStyle = "Absatz-Standardschriftart"
    $s3.Priority = 99
    $s3.LinkStyle = "Kommentartext"
    $s3.Font.Size = 10
    $s3.Font.SizeBi = 10

    $s4 = $d.Styles.Add("Kommentarthema", 1)
    $s4.NameLocal = "annotation subject"
    $s4.BaseStyle = "Kommentartext"
    $s4.NextParagraphStyle = "Kommentartext"
    $s4.Priority = 99
    $s4.UnhideWhenUsed = $true
    $s4.LinkStyle = "KommentarthemaZchn"
    $s4.Font.Bold = $true
    $s4.Font.BoldBi = $true

    $s5 = $d.Styles.Add("KommentarthemaZchn", 2)
    $s5.NameLocal = "Kommentarthema Zchn"
    $s5.BaseStyle = "KommentartextZchn"
    $s5.Priority = 99
    $s5.LinkStyle = "Kommentarthema"
    $s5.Font.Bold = $true
    $s5.Font.BoldBi = $true
    $s5.Font.Size = 10
    $s5.Font.SizeBi = 10

    $s6 = $d.Styles.Add("Sprechblasentext", 1)
    $s6.NameLocal = "Balloon Text"
    $s6.BaseStyle = "Standard"
    $s6.Priority = 99
    $s6.UnhideWhenUsed = $true
    $s6.LinkStyle = "SprechblasentextZchn"
    $s6.ParagraphFormat.SpaceAfter = 0
    $s6.ParagraphFormat.LineSpacingRule = 0
    $s6.Font.NameAscii = "Segoe UI"
    $s6.Font.NameOther = "Segoe UI"
    $s6.Font.NameBi = "Segoe UI"
    $s6.Font.Size = 9
    $s6.Font.SizeBi = 9

    $s7 = $d.Styles.Add("SprechblasentextZchn", 2)
    $s7.NameLocal = "Sprechblasentext Zchn"
    $s7.BaseStyle = "Absatz-Standardschriftart"
    $s7.Priority = 99
    $s7.LinkStyle = "Sprechblasentext"
    $s7.Font.NameAscii = "Segoe UI"
    $s7.Font.NameOther = "Segoe UI"
    $s7.Font.NameBi = "Segoe UI"
    $s7.Font.Size = 9
    $s7.Font.SizeBi = 9
} catch {
    Write-Host "Style gallery registration skipped: $_"
}
